$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 959.1
$ws.Range("I101").Value = 456.83334
$ws.Range("J101").Value = 1712.5
$ws.Range("K101").Value = 1370.50002
$ws.Range("L101").Value = 5137.5
$ws.Range("M101").Value = 251.4999800000001
$ws.Range("N101").Value = -8381.5
$ws.Range("H121").Value = 1094.04
$ws.Range("J121").Value = 1254.762
$ws.Range("L121").Value = 3764.286
$ws.Range("N121").Value = -7258.286
$ws.Range("H138").Value = 4504.846
$ws.Range("I138").Value = 5966.143
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 17898.429
$ws.Range("L138").Value = 8400
$ws.Range("M138").Value = -12758.429
$ws.Range("N138").Value = -18680
$ws.Range("H141").Value = 1611.75
$ws.Range("I141").Value = 984.8570999999999
$ws.Range("K141").Value = 2954.5713
$ws.Range("M141").Value = 2225.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1073.84
$ws.Range("I32").Value = 1015.9785
$ws.Range("K32").Value = 1015.9785
$ws.Range("M32").Value = -728.9785000000001
$ws.Range("H61").Value = 968.125
$ws.Range("I61").Value = 682.5814
$ws.Range("J61").Value = 1912.6154
$ws.Range("K61").Value = 682.5814
$ws.Range("L61").Value = 1912.6154
$ws.Range("M61").Value = -470.5814
$ws.Range("N61").Value = -2336.6154
$ws.Range("H76").Value = 27000
$ws.Range("I76").Value = 34000
$ws.Range("J76").Value = 20000
$ws.Range("K76").Value = 34000
$ws.Range("L76").Value = 20000
$ws.Range("M76").Value = -33662
$ws.Range("N76").Value = -20676
$ws.Range("H79").Value = 27000
$ws.Range("I79").Value = 34000
$ws.Range("J79").Value = 20000
$ws.Range("K79").Value = 34000
$ws.Range("L79").Value = 20000
$ws.Range("M79").Value = -32830
$ws.Range("N79").Value = -22340
$ws.Range("H136").Value = 968.125
$ws.Range("I136").Value = 682.5814
$ws.Range("J136").Value = 1912.6154
$ws.Range("K136").Value = 2047.7442
$ws.Range("L136").Value = 5737.8462
$ws.Range("M136").Value = 502.2557999999999
$ws.Range("N136").Value = -10837.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 713.7692
$ws.Range("I94").Value = 536.86957
$ws.Range("J94").Value = 2070
$ws.Range("K94").Value = 536.86957
$ws.Range("L94").Value = 2070
$ws.Range("M94").Value = -85.86956999999995
$ws.Range("N94").Value = -2972

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4001941
$ws.Range("I31").Value = 2150.4443
$ws.Range("J31").Value = 8697347
$ws.Range("K31").Value = 2150.4443
$ws.Range("L31").Value = 8697347
$ws.Range("M31").Value = -1855.4443
$ws.Range("N31").Value = -8697937
$ws.Range("H34").Value = 4001941
$ws.Range("I34").Value = 2150.4443
$ws.Range("J34").Value = 8697347
$ws.Range("K34").Value = 2150.4443
$ws.Range("L34").Value = 8697347
$ws.Range("M34").Value = -1948.4443
$ws.Range("N34").Value = -8697751
$ws.Range("H58").Value = 749.1707
$ws.Range("I58").Value = 811.4483
$ws.Range("J58").Value = 598.6667
$ws.Range("K58").Value = 811.4483
$ws.Range("L58").Value = 598.6667
$ws.Range("M58").Value = -608.4483
$ws.Range("N58").Value = -1004.6667
$ws.Range("H136").Value = 749.1707
$ws.Range("I136").Value = 811.4483
$ws.Range("J136").Value = 598.6667
$ws.Range("K136").Value = 2434.3449
$ws.Range("L136").Value = 1796.0001
$ws.Range("M136").Value = 115.6550999999999
$ws.Range("N136").Value = -6896.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 66666784
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 100000150
$ws.Range("K26").Value = 150
$ws.Range("L26").Value = 300000450
$ws.Range("M26").Value = 138
$ws.Range("N26").Value = -300001026
$ws.Range("H56").Value = 10041.1
$ws.Range("I56").Value = 10041.1
$ws.Range("K56").Value = 10041.1
$ws.Range("M56").Value = -9511.1
$ws.Range("H68").Value = 1355.303
$ws.Range("J68").Value = 1499
$ws.Range("L68").Value = 4497
$ws.Range("N68").Value = -6119
$ws.Range("H71").Value = 1355.303
$ws.Range("J71").Value = 1499
$ws.Range("L71").Value = 13491
$ws.Range("N71").Value = -21603
$ws.Range("H76").Value = 11352
$ws.Range("I76").Value = 750
$ws.Range("J76").Value = 13119
$ws.Range("K76").Value = 2250
$ws.Range("L76").Value = 39357
$ws.Range("M76").Value = -1867
$ws.Range("N76").Value = -40123
$ws.Range("H79").Value = 11352
$ws.Range("I79").Value = 750
$ws.Range("J79").Value = 13119
$ws.Range("K79").Value = 2250
$ws.Range("L79").Value = 39357
$ws.Range("M79").Value = -924
$ws.Range("N79").Value = -42009
$ws.Range("H122").Value = 868.56525
$ws.Range("J122").Value = 842.3461
$ws.Range("L122").Value = 7581.1149
$ws.Range("N122").Value = -12481.1149
$ws.Range("H134").Value = 1571.6111
$ws.Range("I134").Value = 1406.8462
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 4220.5386
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 849.4614000000001
$ws.Range("N134").Value = -16140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 56325.11
$ws.Range("I132").Value = 101770.15
$ws.Range("J132").Value = 2860.353
$ws.Range("K132").Value = 305310.45
$ws.Range("L132").Value = 8581.059000000001
$ws.Range("M132").Value = -302780.45
$ws.Range("N132").Value = -13641.059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 29000
$ws.Range("I62").Value = 29000
$ws.Range("K62").Value = 29000
$ws.Range("M62").Value = -28376
$ws.Range("H65").Value = 29000
$ws.Range("I65").Value = 29000
$ws.Range("K65").Value = 87000
$ws.Range("M65").Value = -83880
$ws.Range("H76").Value = 14288
$ws.Range("I76").Value = 14288
$ws.Range("K76").Value = 14288
$ws.Range("M76").Value = -13950
$ws.Range("H79").Value = 14288
$ws.Range("I79").Value = 14288
$ws.Range("K79").Value = 14288
$ws.Range("M79").Value = -13118
$ws.Range("H93").Value = 1352539.9
$ws.Range("I93").Value = 1931570.9
$ws.Range("J93").Value = 1467.6666
$ws.Range("K93").Value = 1931570.9
$ws.Range("L93").Value = 1467.6666
$ws.Range("M93").Value = -1930322.9
$ws.Range("N93").Value = -3963.6666
$ws.Range("H132").Value = 7015.1797
$ws.Range("I132").Value = 11285.1
$ws.Range("J132").Value = 2520.5264
$ws.Range("K132").Value = 33855.3
$ws.Range("L132").Value = 7561.5792
$ws.Range("M132").Value = -31325.3
$ws.Range("N132").Value = -12621.5792
$ws.Range("H136").Value = 7782.684
$ws.Range("I136").Value = 11418.363
$ws.Range("K136").Value = 34255.089
$ws.Range("M136").Value = -31705.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13056.2
$ws.Range("I136").Value = 41920.668
$ws.Range("J136").Value = 685.7143
$ws.Range("K136").Value = 125762.004
$ws.Range("L136").Value = 2057.1429
$ws.Range("M136").Value = -123212.004
$ws.Range("N136").Value = -7157.1429
